$d = $word.ActiveDocument

# --- Edit 2 first: drop the original "_GoBack" bookmark that sits between
#     the "...SourceTree" run and the trailing "  " run, and merge that
#     text into a single run "...SourceTree  ". (Done before edit 1 since
#     a document can only have one bookmark named "_GoBack" at a time.) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$targetText = "because of issues with the AGL proxy, GitExtensions is preferred respect to SourceTree  "
$find = $d.Content
$found = $find.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $find.Find.Found) {
    # Fallback: the two trailing spaces may not have matched as part of
    # the search string; locate the sentence alone and grow the range to
    # include whatever immediately follows it up to the paragraph end.
    $find = $d.Content
    $found = $find.Find.Execute("because of issues with the AGL proxy, GitExtensions is preferred respect to SourceTree", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $para = $find.Paragraphs(1)
    $find.End = $para.Range.End - 1
}

$mergeRange = $d.Range($find.Start, $find.End)
$mergeStart = $mergeRange.Start
$mergeRange.Delete()

$anchor = $d.Range($mergeStart, $mergeStart)
$anchor.InsertAfter($targetText)

# --- Edit 1: insert "TEST_A " run + a new "_GoBack" bookmark right before
#     the document's opening "LIST OF SOFTWARE..." heading text. ---
$headingPara = $d.Paragraphs(1)
$insertPoint = $headingPara.Range.Duplicate
$insertPoint.Collapse(1)
$insertPoint.InsertBefore("TEST_A ")

# New run is exactly 7 characters ("TEST_A " incl. trailing space) so the
# bookmark collapses right after it, before the existing heading text.
$goBackRange = $d.Range(7, 7)
$null = $d.Bookmarks.Add("_GoBack", $goBackRange)
